# Update the "dSF" (column F) values to reflect re-pulled data / mean calc.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F3"  = -6
    "F4"  = 0
    "F6"  = -1
    "F10" = -8
    "F11" = -6
    "F12" = -7
    "F14" = 0
    "F15" = 2
    "F18" = -5
    "F19" = -6
    "F25" = 3
    "F29" = -9
    "F30" = -6
    "F31" = -4
    "F34" = -3
    "F36" = -6
    "F38" = 2
    "F40" = -1
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
